$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.676.75"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "3.591.32"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'610.07"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'148.68"
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("D7").Value = "3.590.69"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "'0.417"
$ws.Range("D13").Value = "4.199.07"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "'0.0000210"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'30.20"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "3.583.43"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "66.752.20"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'11.48"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").Value = "'6.34"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "'15.19"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "'432.20"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'0.627"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").Value = "'79.12"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Value = "3.736.74"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'0.0000121"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "'8.25"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").Value = "'9.34"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "3.586.29"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.46"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'25.54"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("E35").Value = "  -3.32%  "
$ws.Range("D36").Value = "'7.89"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("D39").Value = "'5.65"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'173.62"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("D41").Value = "'0.0859"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "'46.14"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "'2.57"
$ws.Range("E46").Value = "  +7.34%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").Value = "'25.09"
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("D50").Value = "'23.94"
$ws.Range("E50").Value = "  +4.36%  "
$ws.Range("E51").Value = "  +1.17%  "
